$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.617.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.18%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.766.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.91%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.67%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'598.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.22%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'162.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.70%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.763.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.34%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.35%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'6.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.11%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.63%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'35.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.56%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.397.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.95%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.763.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.46%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.605.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.27%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'18.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.42%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'455.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.18%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -4.69%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.689"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.26%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'82.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.79%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0000142"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -6.40%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'11.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.81%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.00%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.97%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.915.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.74%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.54%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.45%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.64%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.52%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'8.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.31%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0987"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.27%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.22%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.40%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B41").Value = "'dogwifhat"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'3.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -6.51%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'FirstDigitalUSD"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.07%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D44").Value = "'43.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.28%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'47.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'151.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.77%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.294"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.25%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.95%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.76%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'383.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.81%  "
$ws.Range("E51").Style = "Normal"
